# Apply Betfair Back/Lay odds updates for 2025-11-20 (diff-driven cell updates)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 2.52
$ws.Range("G2").Value = 2.54
$ws.Range("H2").Value = 4.3
$ws.Range("I2").Value = 4.4
$ws.Range("J2").Value = 2.68
$ws.Range("K2").Value = 2.7
$ws.Range("L2").Value = 1.67
$ws.Range("N2").Value = 2.5
$ws.Range("O2").Value = 1.65
$ws.Range("P2").Value = 1.45
$ws.Range("Q2").Value = 3
$ws.Range("V2").Value = 1.29
$ws.Range("W2").Value = 1.65
$ws.Range("X2").Value = 7
$ws.Range("Z2").Value = 32
$ws.Range("AA2").Value = 900
$ws.Range("AB2").Value = 7.2
$ws.Range("AD2").Value = 19
$ws.Range("AE2").Value = 1000
$ws.Range("AF2").Value = 13
$ws.Range("AG2").Value = 12
$ws.Range("AH2").Value = 28
$ws.Range("AJ2").Value = 38
$ws.Range("AN2").Value = 44

# Row 3
$ws.Range("F3").Value = 3.35
$ws.Range("G3").Value = 3.5
$ws.Range("H3").Value = 2.28
$ws.Range("I3").Value = 2.34
$ws.Range("P3").Value = 2.18
$ws.Range("Q3").Value = 1.79
$ws.Range("R3").Value = 1.47
$ws.Range("S3").Value = 3
$ws.Range("T3").Value = 1.66
$ws.Range("U3").Value = 2.34
$ws.Range("V3").Value = 1.74
$ws.Range("W3").Value = 1.4
$ws.Range("X3").Value = 18.5
$ws.Range("AB3").Value = 16
$ws.Range("AC3").Value = 9
$ws.Range("AE3").Value = 24
$ws.Range("AF3").Value = 25
$ws.Range("AI3").Value = 34
$ws.Range("AJ3").Value = 220
$ws.Range("AK3").Value = 95
$ws.Range("AM3").Value = 200
$ws.Range("AO3").Value = 14.5

# Row 4
$ws.Range("F4").Value = 2.44
$ws.Range("G4").Value = 2.56
$ws.Range("H4").Value = 3
$ws.Range("J4").Value = 3.4
$ws.Range("P4").Value = 1.81
$ws.Range("R4").Value = 1.31
$ws.Range("V4").Value = 1.45
$ws.Range("W4").Value = 1.64
$ws.Range("AA4").Value = 270
$ws.Range("AE4").Value = 130
$ws.Range("AF4").Value = 15.5
$ws.Range("AJ4").Value = 36
$ws.Range("AK4").Value = 29
$ws.Range("AL4").Value = 130
$ws.Range("AN4").Value = 26
$ws.Range("AO4").Value = 970

# Row 5
$ws.Range("F5").Value = 1.9
$ws.Range("G5").Value = 1.98
$ws.Range("H5").Value = 5.6
$ws.Range("I5").Value = 6.6
$ws.Range("K5").Value = 3.35
$ws.Range("P5").Value = 1.46
$ws.Range("T5").Value = 2.36
$ws.Range("U5").Value = 1.6
$ws.Range("V5").Value = 1.2
$ws.Range("W5").Value = 2.02
$ws.Range("Z5").Value = 48
$ws.Range("AA5").Value = 200
$ws.Range("AB5").Value = 5.9
$ws.Range("AD5").Value = 26
$ws.Range("AE5").Value = 140
$ws.Range("AF5").Value = 10
$ws.Range("AH5").Value = 55
$ws.Range("AJ5").Value = 24
$ws.Range("AN5").Value = 28
$ws.Range("AO5").Value = 270

# Row 6
$ws.Range("F6").Value = 4.2
$ws.Range("H6").Value = 2.12
$ws.Range("I6").Value = 2.16
$ws.Range("J6").Value = 3.35
$ws.Range("K6").Value = 3.4
$ws.Range("O6").Value = 1.5
$ws.Range("Q6").Value = 2.52
$ws.Range("U6").Value = 1.82
$ws.Range("V6").Value = 1.86
$ws.Range("W6").Value = 1.3
$ws.Range("X6").Value = 10
$ws.Range("AN6").Value = 340
$ws.Range("AO6").Value = 24

# Row 7
$ws.Range("Q7").Value = 2.08

# Row 8
$ws.Range("F8").Value = 2.48
$ws.Range("G8").Value = 2.64
$ws.Range("I8").Value = 3.5
$ws.Range("J8").Value = 3.15
$ws.Range("K8").Value = 3.35
$ws.Range("N8").Value = 3.6
$ws.Range("P8").Value = 1.86
$ws.Range("Q8").Value = 2.06
$ws.Range("T8").Value = 1.78
$ws.Range("W8").Value = 1.6
$ws.Range("Y8").Value = 42
$ws.Range("AA8").Value = 900
$ws.Range("AB8").Value = 15
$ws.Range("AG8").Value = 12
$ws.Range("AJ8").Value = 900
$ws.Range("AK8").Value = 70
$ws.Range("AL8").Value = 110
$ws.Range("AO8").Value = 980

# Row 9
$ws.Range("G9").Value = 1.55
$ws.Range("J9").Value = 4.6
$ws.Range("K9").Value = 4.7
$ws.Range("Q9").Value = 1.72
$ws.Range("S9").Value = 2.88
$ws.Range("U9").Value = 2.1
$ws.Range("W9").Value = 2.78
$ws.Range("AA9").Value = 250
$ws.Range("AI9").Value = 85
$ws.Range("AO9").Value = 100

# Row 10
$ws.Range("S10").Value = 3.65
$ws.Range("V10").Value = 1.23

# Row 11
$ws.Range("F11").Value = 2.5
$ws.Range("G11").Value = 2.52
$ws.Range("J11").Value = 2.98
$ws.Range("K11").Value = 3
$ws.Range("N11").Value = 2.66
$ws.Range("O11").Value = 1.59
$ws.Range("W11").Value = 1.65
$ws.Range("AA11").Value = 80
$ws.Range("AN11").Value = 40

# Row 12
$ws.Range("F12").Value = 1.99

# Row 13
$ws.Range("N13").Value = 2.76
$ws.Range("Q13").Value = 2.66

# Row 14
$ws.Range("H14").Value = 5.8
$ws.Range("Q14").Value = 1.93
$ws.Range("T14").Value = 1.89
$ws.Range("Y14").Value = 25
$ws.Range("AK14").Value = 22

# Row 16
$ws.Range("H16").Value = 3.65
$ws.Range("I16").Value = 3.8
